# Remove the now-redundant "reviews_count" column (column E).
# All of its cells are empty, so this is a straightforward column delete
# that shifts reviews_average/latitude/longitude/is_permanently_closed/
# gmaps_link/latest_review_date (columns F:K) one position to the left
# (E:J), and shrinks the used range from A1:K65 to A1:J65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").EntireColumn.Delete()
